$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from B1 (bold, bordered, centered) onto the
# new H1 header cell, then set its text to "Label".
$ws.Range("B1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Label"

# New "Label" column: 0 for Control patients, 1 for MDD patients.
# The worksheet repeats the same 10 patients twice (rows 2-11 and 12-21).
$labels = @(0, 0, 0, 0, 0, 1, 1, 1, 1, 1)

for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(2 + $i, 8).Value = $labels[$i]
    $ws.Cells.Item(12 + $i, 8).Value = $labels[$i]
}
